$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Total" column header (N3) ---
$ws.Range("N3").Value = "Total"

# --- Add Total column formulas.
#     N4 is its own (non-shared) formula; N5:N9 were filled together as one
#     shared-formula group (relative SUM one row down each time). ---
$ws.Range("N4").Formula = "=SUM(B4:M4)"
$ws.Range("N5:N9").Formula = "=SUM(B5:M5)"

# --- Add Total row label (A9) ---
$ws.Range("A9").Value = "Total"

# --- Add Total row formulas. The original author filled several
#     contiguous blocks separately (each becomes its own shared-formula
#     group), rather than the whole row in one shot. ---
$ws.Range("B9").Formula = "=SUM(B4:B8)"
$ws.Range("C9:I9").Formula = "=SUM(C4:C8)"
$ws.Range("J9").Formula = "=SUM(J4:J8)"
$ws.Range("K9").Formula = "=SUM(K4:K8)"
$ws.Range("L9").Formula = "=SUM(L4:L8)"
$ws.Range("M9:N9").Formula = "=SUM(M4:M8)"

# --- Scroll the view down a bit and land the selection on I5, matching
#     the author's final cursor position when they saved. ---
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("I5").Select()
